# Weekly update: insert two new data rows (row 369 and 370) at the top of the
# "Betarraga" data block (which starts at row 369), pushing the existing
# rows 369:435 down to 371:437.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 new rows before row 369; existing data (rows 369-435) shifts down
# to rows 371-437.
$ws.Rows("369:370").Insert()

# --- Row 369: new "Primera" quality data point ---
$ws.Cells.Item(369, 1).Value = 7
$ws.Cells.Item(369, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(369, 3).Value = "Ñuble"
$ws.Cells.Item(369, 4).Value = 44816
$ws.Cells.Item(369, 5).Value = 16
$ws.Cells.Item(369, 6).Value = 100114014
$ws.Cells.Item(369, 7).Value = "Betarraga"
$ws.Cells.Item(369, 8).Value = "Sin especificar"
$ws.Cells.Item(369, 9).Value = "Primera"
$ws.Cells.Item(369, 10).Value = 300
$ws.Cells.Item(369, 11).Value = 900
$ws.Cells.Item(369, 12).Value = 1000
$ws.Cells.Item(369, 13).Value = 950
$ws.Cells.Item(369, 14).Value = '$/paquete 5 unidades'
$ws.Cells.Item(369, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(369, 16).Value = 190
$ws.Cells.Item(369, 17).Value = 5
$ws.Cells.Item(369, 18).Value = "Hortaliza"

# --- Row 370: new "Segunda" quality data point ---
$ws.Cells.Item(370, 1).Value = 7
$ws.Cells.Item(370, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(370, 3).Value = "Ñuble"
$ws.Cells.Item(370, 4).Value = 44816
$ws.Cells.Item(370, 5).Value = 16
$ws.Cells.Item(370, 6).Value = 100114014
$ws.Cells.Item(370, 7).Value = "Betarraga"
$ws.Cells.Item(370, 8).Value = "Sin especificar"
$ws.Cells.Item(370, 9).Value = "Segunda"
$ws.Cells.Item(370, 10).Value = 240
$ws.Cells.Item(370, 11).Value = 700
$ws.Cells.Item(370, 12).Value = 800
$ws.Cells.Item(370, 13).Value = 750
$ws.Cells.Item(370, 14).Value = '$/paquete 5 unidades'
$ws.Cells.Item(370, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(370, 16).Value = 150
$ws.Cells.Item(370, 17).Value = 5
$ws.Cells.Item(370, 18).Value = "Hortaliza"

# Apply the same date style (numFmtId 165 date/time) that the rest of
# column D uses, matching the surrounding cells (e.g. D371).
$ws.Range("D369:D370").NumberFormat = $ws.Range("D371").NumberFormat
